$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new history row (row 16), matching the formatting of row 15 above it
$ws.Range("A15:B15").Copy()
$ws.Range("A16:B16").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A16").Value = 45750
$ws.Range("B16").Value = 0.40972222222222221
$ws.Range("C16").Value = "Futconnect0304 0950"
$ws.Range("D16").Value = "Responsividade das páginas."

$ws.Range("D17").Select()
